# Generate Report for Handoff
# Adds a new row (for file 3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6) to each of
# the three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$commitSha = "b5970467c1098d4409e2b37952e381f3b98f6e23"
$newFile = "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md"
$newFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newFile"
$newFileDisplay = "e2e\3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(3, 1).Value() = $newFile
$wsOverview.Cells.Item(3, 2).Value() = $newFileDisplay
$wsOverview.Cells.Item(3, 3).Value() = ".md"
$wsOverview.Cells.Item(3, 4).Value() = "'"
$wsOverview.Cells.Item(3, 5).Value() = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value() = "Ready for handoff"
$wsOverview.Cells.Item(3, 7).Value() = "2016-08-28 02:39:48"
$wsOverview.Cells.Item(3, 7).NumberFormat() = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newFileUrl, [Type]::Missing, [Type]::Missing, $newFileDisplay) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Cells.Item(3, 1).Value() = $newFile
$wsZhCn.Cells.Item(3, 2).Value() = ".md"
$wsZhCn.Cells.Item(3, 3).Value() = "Ready for handoff"
$wsZhCn.Cells.Item(3, 4).Value() = "e2e"
$wsZhCn.Cells.Item(3, 5).Value() = "ht"
$wsZhCn.Cells.Item(3, 6).Value() = "'False"
$wsZhCn.Cells.Item(3, 7).Value() = "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.fff2734d607640bd36765059c09fb28d3bc65cc1.zh-cn.xlf"
$wsZhCn.Cells.Item(3, 8).Value() = "2016-08-28 02:39:43"
$wsZhCn.Cells.Item(3, 8).NumberFormat() = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3, 9).Value() = "'"
$wsZhCn.Cells.Item(3, 10).Value() = "'"
$wsZhCn.Cells.Item(3, 11).Value() = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(3, 11).NumberFormat() = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3, 12).Value() = "'"
$wsZhCn.Cells.Item(3, 13).Value() = "'True"
$wsZhCn.Cells.Item(3, 14).Value() = "'"
$wsZhCn.Cells.Item(3, 15).Value() = "'False"
$wsZhCn.Cells.Item(3, 16).Value() = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newFileUrl, [Type]::Missing, [Type]::Missing, $newFile) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Cells.Item(3, 1).Value() = $newFile
$wsDeDe.Cells.Item(3, 2).Value() = ".md"
$wsDeDe.Cells.Item(3, 3).Value() = "Ready for handoff"
$wsDeDe.Cells.Item(3, 4).Value() = "e2e"
$wsDeDe.Cells.Item(3, 5).Value() = "ht"
$wsDeDe.Cells.Item(3, 6).Value() = "'False"
$wsDeDe.Cells.Item(3, 7).Value() = "3e537be1-bbc8-4711-8dfd-c0e1ac23b1c6.fff2734d607640bd36765059c09fb28d3bc65cc1.de-de.xlf"
$wsDeDe.Cells.Item(3, 8).Value() = "2016-08-28 02:39:48"
$wsDeDe.Cells.Item(3, 8).NumberFormat() = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3, 9).Value() = "'"
$wsDeDe.Cells.Item(3, 10).Value() = "'"
$wsDeDe.Cells.Item(3, 11).Value() = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(3, 11).NumberFormat() = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3, 12).Value() = "'"
$wsDeDe.Cells.Item(3, 13).Value() = "'True"
$wsDeDe.Cells.Item(3, 14).Value() = "'"
$wsDeDe.Cells.Item(3, 15).Value() = "'False"
$wsDeDe.Cells.Item(3, 16).Value() = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newFileUrl, [Type]::Missing, [Type]::Missing, $newFile) | Out-Null

Write-Host "Report row added for $newFile across Overview, zh-cn, de-de sheets."
